# Updated cryptos list (Price / Volume(1h) columns) with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume(1h) (E) columns to be treated as text before
# writing the new values, so strings like "1.005" or "24.561.39" are stored as
# literal text (matching the original inline-string cell type) instead of being
# auto-converted into numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.561.39"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.701.27"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "307.89"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.3730"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "48.88"
$ws.Range("E8").Value = "  +2.73%  "
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Value = "0.07423"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "20.76"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "6.202"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "6.889"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "1.702.07"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "0.06684"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "82.95"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").Value = "6.316"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").Value = "13.10"
$ws.Range("E23").Value = "  +8.63%  "
$ws.Range("D24").Value = "24.545.55"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "2.420"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "2.756"
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("D27").Value = "20.03"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").Value = "149.08"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").Value = "130.70"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("D30").Value = "1.889.50"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").Value = "1.171"
$ws.Range("E31").Value = "  +17.84%  "
$ws.Range("D32").Value = "6.675"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").Value = "4.212"
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").Value = "0.08746"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("D35").Value = "1.767"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").Value = "13.48"
$ws.Range("E36").Value = "  +7.08%  "
$ws.Range("D37").Value = "5.475"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "0.06489"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "8.893"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "0.02362"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "0.2208"
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("D42").Value = "1.269"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "0.6368"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "13.74"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("D46").Value = "0.6047"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").Value = "3.795"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "2.101"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "0.07228"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "78.65"
$ws.Range("E51").Value = "  +2.47%  "

# Clear the temporary text formatting so the cells keep the same (absent)
# style index as before the edit.
$dataRange.ClearFormats()
